$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "{0}{1}" -f [char]39, $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "69.420.46"
Set-TextValue "E2" "  -2.33%  "

Set-TextValue "D3" "3.694.92"
Set-TextValue "E3" "  -2.98%  "

Set-TextValue "E4" "  -0.03%  "

Set-TextValue "D5" "691.18"
Set-TextValue "E5" "  -1.60%  "

Set-TextValue "D6" "162.36"
Set-TextValue "E6" "  -5.45%  "

Set-TextValue "D7" "3.693.74"
Set-TextValue "E7" "  -2.96%  "

Set-TextValue "E8" "  +0.03%  "

Set-TextValue "E9" "  -4.81%  "

Set-TextValue "E10" "  -8.19%  "

Set-TextValue "D11" "7.37"
Set-TextValue "E11" "  -1.99%  "

Set-TextValue "E12" "  -5.29%  "

Set-TextValue "D13" "0.0000238"
Set-TextValue "E13" "  -5.18%  "

Set-TextValue "D14" "33.36"
Set-TextValue "E14" "  -7.11%  "

Set-TextValue "D15" "4.316.53"
Set-TextValue "E15" "  -3.02%  "

Set-TextValue "D16" "3.694.77"
Set-TextValue "E16" "  -2.71%  "

Set-TextValue "D17" "69.444.77"
Set-TextValue "E17" "  -2.38%  "

Set-TextValue "D18" "0.114"
Set-TextValue "E18" "  -0.66%  "

Set-TextValue "D19" "16.17"
Set-TextValue "E19" "  -7.56%  "

Set-TextValue "E20" "  -7.89%  "

Set-TextValue "D21" "480.27"
Set-TextValue "E21" "  -6.45%  "

Set-TextValue "D22" "9.99"
Set-TextValue "E22" "  -5.89%  "

Set-TextValue "E23" "  -7.04%  "

Set-TextValue "D24" "79.83"
Set-TextValue "E24" "  -4.79%  "

Set-TextValue "D25" "3.840.48"
Set-TextValue "E25" "  -2.97%  "

Set-TextValue "E26" "  -8.93%  "

Set-TextValue "E27" "  +0.09%  "

Set-TextValue "D28" "11.36"
Set-TextValue "E28" "  -5.69%  "

Set-TextValue "D29" "9.51"
Set-TextValue "E29" "  -8.39%  "

Set-TextValue "E30" "  -10.11%  "

Set-TextValue "E31" "  -10.12%  "

Set-TextValue "D32" "6.83"
Set-TextValue "E32" "  -7.57%  "

Set-TextValue "E33" "  -7.59%  "

Set-TextValue "E34" "  -5.71%  "

Set-TextValue "D35" "0.999"
Set-TextValue "E35" "  +0.44%  "

Set-TextValue "D36" "26.95"
Set-TextValue "E36" "  -7.09%  "

Set-TextValue "D37" "3.665.50"
Set-TextValue "E37" "  -2.76%  "

Set-TextValue "D38" "8.47"
Set-TextValue "E38" "  -7.22%  "

Set-TextValue "D39" "6.31"
Set-TextValue "E39" "  +5.43%  "

Set-TextValue "D40" "2.34"
Set-TextValue "E40" "  -2.18%  "

Set-TextValue "E41" "  -7.90%  "

Set-TextValue "E43" "  -0.02%  "

Set-TextValue "D44" "0.953"
Set-TextValue "E44" "  -6.47%  "

Set-TextValue "D45" "163.66"
Set-TextValue "E45" "  -5.30%  "

Set-TextValue "D46" "48.09"
Set-TextValue "E46" "  -2.52%  "

Set-TextValue "D47" "30.10"
Set-TextValue "E47" "  +3.33%  "

Set-TextValue "D48" "2.80"
Set-TextValue "E48" "  -15.12%  "

Set-TextValue "D49" "1.15"
Set-TextValue "E49" "  -0.60%  "

Set-TextValue "D50" "1.35"
Set-TextValue "E50" "  -2.06%  "

Set-TextValue "E51" "  -9.00%  "

